# Fruta / hortaliza, semanal
# Insert a new weekly price-report row at the top of the data block
# (row 201), pushing the existing rows 201-259 down to 202-260.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 201..259 down to 202..260, leaving a blank row 201 to fill in.
$ws.Rows.Item(201).Insert()

# Populate the new row 201 with the latest observation.
$ws.Range("A201").Value = 3
$ws.Range("B201").Value = "Femacal de La Calera"
$ws.Range("C201").Value = "Coquimbo"
$ws.Range("D201").Value = 44489
$ws.Range("E201").Value = 5
$ws.Range("F201").Value = 100112028
$ws.Range("G201").Value = "Sandia"
$ws.Range("H201").Value = "Sin especificar"
$ws.Range("I201").Value = "Primera"
$ws.Range("J201").Value = 220
$ws.Range("K201").Value = 800
$ws.Range("L201").Value = 800
$ws.Range("M201").Value = 800
$ws.Range("N201").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O201").Value = "Perú"
$ws.Range("P201").Value = 800
$ws.Range("Q201").Value = 1
$ws.Range("R201").Value = "Hortaliza"
